$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E4").Value = "2016-03-19 00:34:18"
$wsZhCn.Range("H4").Value = "2016-03-19 00:34:37"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E4").Value = "2016-03-19 00:34:22"
$wsDeDe.Range("H4").Value = "2016-03-19 00:34:42"
